# Auto-generated PowerShell Excel COM script
# Applies the diff: renumber/insert '2022-Q4' sheet + update summary sheet

$wb = $excel.ActiveWorkbook

# ================= Step 1: Update the summary sheet (总计) =================
$summary = $wb.Worksheets.Item(1)

# Insert a new row at position 2 for the 2022-Q4 entry; existing rows shift down
$summary.Rows(2).Insert()
$summary.Range('B2:D2').ClearFormats()

$summary.Range('A2').Value = 0
$summary.Range('A2').Font.Bold = $true
$summary.Range('A2').HorizontalAlignment = -4108
$summary.Range('A2').VerticalAlignment = -4160
$summary.Range('A2').Borders.LineStyle = 1

$summary.Range('B2').Value = "2022-Q4"
$summary.Range('C2').Value = 42
$summary.Range('D2').Value = 8.15

# Re-index the rows that shifted down (0-based index in column A)
$summary.Range('A3').Value = 1
$summary.Range('A4').Value = 2
$summary.Range('A5').Value = 3
$summary.Range('A6').Value = 4

# ================= Step 2: Insert the new '2022-Q4' worksheet =================
# Position it right after the summary sheet (总计), before 2022-Q3
$afterSheet = $wb.Worksheets.Item(1)
$q4 = $wb.Worksheets.Add($null, $afterSheet)
$q4.Name = '2022-Q4'

# ---- Header row ----
$q4.Range('B1').Value = "基金代码"
$q4.Range('B1').Font.Bold = $true
$q4.Range('B1').HorizontalAlignment = -4108
$q4.Range('B1').VerticalAlignment = -4160
$q4.Range('B1').Borders.LineStyle = 1
$q4.Range('C1').Value = "基金名称"
$q4.Range('C1').Font.Bold = $true
$q4.Range('C1').HorizontalAlignment = -4108
$q4.Range('C1').VerticalAlignment = -4160
$q4.Range('C1').Borders.LineStyle = 1
$q4.Range('D1').Value = "基金规模"
$q4.Range('D1').Font.Bold = $true
$q4.Range('D1').HorizontalAlignment = -4108
$q4.Range('D1').VerticalAlignment = -4160
$q4.Range('D1').Borders.LineStyle = 1
$q4.Range('E1').Value = "股票总仓位"
$q4.Range('E1').Font.Bold = $true
$q4.Range('E1').HorizontalAlignment = -4108
$q4.Range('E1').VerticalAlignment = -4160
$q4.Range('E1').Borders.LineStyle = 1
$q4.Range('F1').Value = "仓位占比"
$q4.Range('F1').Font.Bold = $true
$q4.Range('F1').HorizontalAlignment = -4108
$q4.Range('F1').VerticalAlignment = -4160
$q4.Range('F1').Borders.LineStyle = 1
$q4.Range('G1').Value = "持有市值(亿元)"
$q4.Range('G1').Font.Bold = $true
$q4.Range('G1').HorizontalAlignment = -4108
$q4.Range('G1').VerticalAlignment = -4160
$q4.Range('G1').Borders.LineStyle = 1
$q4.Range('H1').Value = "仓位排名"
$q4.Range('H1').Font.Bold = $true
$q4.Range('H1').HorizontalAlignment = -4108
$q4.Range('H1').VerticalAlignment = -4160
$q4.Range('H1').Borders.LineStyle = 1

# ---- Data rows (index column A is bold/bordered; fund data in B:H) ----
# Row 2: fund 166002
$q4.Range('A2').Value = 0
$q4.Range('A2').Font.Bold = $true
$q4.Range('A2').HorizontalAlignment = -4108
$q4.Range('A2').VerticalAlignment = -4160
$q4.Range('A2').Borders.LineStyle = 1
$q4.Range('B2').NumberFormat = '@'
$q4.Range('B2').Value = "166002"
$q4.Range('C2').NumberFormat = '@'
$q4.Range('C2').Value = "中欧新蓝筹混合 -A"
$q4.Range('D2').NumberFormat = '@'
$q4.Range('D2').Value = "109.21"
$q4.Range('E2').NumberFormat = '@'
$q4.Range('E2').Value = "78.23"
$q4.Range('F2').NumberFormat = '@'
$q4.Range('F2').Value = "2.02"
$q4.Range('G2').NumberFormat = '@'
$q4.Range('G2').Value = "2.2060"
$q4.Range('H2').Value = 5

# Row 3: fund 003634
$q4.Range('A3').Value = 1
$q4.Range('A3').Font.Bold = $true
$q4.Range('A3').HorizontalAlignment = -4108
$q4.Range('A3').VerticalAlignment = -4160
$q4.Range('A3').Borders.LineStyle = 1
$q4.Range('B3').NumberFormat = '@'
$q4.Range('B3').Value = "003634"
$q4.Range('C3').NumberFormat = '@'
$q4.Range('C3').Value = "嘉实农业产业股票A"
$q4.Range('D3').NumberFormat = '@'
$q4.Range('D3').Value = "23.71"
$q4.Range('E3').NumberFormat = '@'
$q4.Range('E3').Value = "92.18"
$q4.Range('F3').NumberFormat = '@'
$q4.Range('F3').Value = "5.72"
$q4.Range('G3').NumberFormat = '@'
$q4.Range('G3').Value = "1.3562"
$q4.Range('H3').Value = 8

# Row 4: fund 012647
$q4.Range('A4').Value = 2
$q4.Range('A4').Font.Bold = $true
$q4.Range('A4').HorizontalAlignment = -4108
$q4.Range('A4').VerticalAlignment = -4160
$q4.Range('A4').Borders.LineStyle = 1
$q4.Range('B4').NumberFormat = '@'
$q4.Range('B4').Value = "012647"
$q4.Range('C4').NumberFormat = '@'
$q4.Range('C4').Value = "中欧洞见一年持有混合"
$q4.Range('D4').NumberFormat = '@'
$q4.Range('D4').Value = "30.44"
$q4.Range('E4').NumberFormat = '@'
$q4.Range('E4').Value = "93.89"
$q4.Range('F4').NumberFormat = '@'
$q4.Range('F4').Value = "2.24"
$q4.Range('G4').NumberFormat = '@'
$q4.Range('G4').Value = "0.6819"
$q4.Range('H4').Value = 10

# Row 5: fund 005106
$q4.Range('A5').Value = 3
$q4.Range('A5').Font.Bold = $true
$q4.Range('A5').HorizontalAlignment = -4108
$q4.Range('A5').VerticalAlignment = -4160
$q4.Range('A5').Borders.LineStyle = 1
$q4.Range('B5').NumberFormat = '@'
$q4.Range('B5').Value = "005106"
$q4.Range('C5').NumberFormat = '@'
$q4.Range('C5').Value = "银华农业产业股票A"
$q4.Range('D5').NumberFormat = '@'
$q4.Range('D5').Value = "12.40"
$q4.Range('E5').NumberFormat = '@'
$q4.Range('E5').Value = "93.03"
$q4.Range('F5').NumberFormat = '@'
$q4.Range('F5').Value = "5.33"
$q4.Range('G5').NumberFormat = '@'
$q4.Range('G5').Value = "0.6609"
$q4.Range('H5').Value = 7

# Row 6: fund 121005
$q4.Range('A6').Value = 4
$q4.Range('A6').Font.Bold = $true
$q4.Range('A6').HorizontalAlignment = -4108
$q4.Range('A6').VerticalAlignment = -4160
$q4.Range('A6').Borders.LineStyle = 1
$q4.Range('B6').NumberFormat = '@'
$q4.Range('B6').Value = "121005"
$q4.Range('C6').NumberFormat = '@'
$q4.Range('C6').Value = "国投瑞银创新动力混合"
$q4.Range('D6').NumberFormat = '@'
$q4.Range('D6').Value = "10.75"
$q4.Range('E6').NumberFormat = '@'
$q4.Range('E6').Value = "89.92"
$q4.Range('F6').NumberFormat = '@'
$q4.Range('F6').Value = "3.10"
$q4.Range('G6').NumberFormat = '@'
$q4.Range('G6').Value = "0.3332"
$q4.Range('H6').Value = 10

# Row 7: fund 015468
$q4.Range('A7').Value = 5
$q4.Range('A7').Font.Bold = $true
$q4.Range('A7').HorizontalAlignment = -4108
$q4.Range('A7').VerticalAlignment = -4160
$q4.Range('A7').Borders.LineStyle = 1
$q4.Range('B7').NumberFormat = '@'
$q4.Range('B7').Value = "015468"
$q4.Range('C7').NumberFormat = '@'
$q4.Range('C7').Value = "嘉实农业产业股票C"
$q4.Range('D7').NumberFormat = '@'
$q4.Range('D7').Value = "5.82"
$q4.Range('E7').NumberFormat = '@'
$q4.Range('E7').Value = "92.18"
$q4.Range('F7').NumberFormat = '@'
$q4.Range('F7').Value = "5.72"
$q4.Range('G7').NumberFormat = '@'
$q4.Range('G7').Value = "0.3329"
$q4.Range('H7').Value = 8

# Row 8: fund 164403
$q4.Range('A8').Value = 6
$q4.Range('A8').Font.Bold = $true
$q4.Range('A8').HorizontalAlignment = -4108
$q4.Range('A8').VerticalAlignment = -4160
$q4.Range('A8').Borders.LineStyle = 1
$q4.Range('B8').NumberFormat = '@'
$q4.Range('B8').Value = "164403"
$q4.Range('C8').NumberFormat = '@'
$q4.Range('C8').Value = "前海开源沪港深农业混合（LOF）A"
$q4.Range('D8').NumberFormat = '@'
$q4.Range('D8').Value = "4.16"
$q4.Range('E8').NumberFormat = '@'
$q4.Range('E8').Value = "88.37"
$q4.Range('F8').NumberFormat = '@'
$q4.Range('F8').Value = "7.77"
$q4.Range('G8').NumberFormat = '@'
$q4.Range('G8').Value = "0.3232"
$q4.Range('H8').Value = 1

# Row 9: fund 501209
$q4.Range('A9').Value = 7
$q4.Range('A9').Font.Bold = $true
$q4.Range('A9').HorizontalAlignment = -4108
$q4.Range('A9').VerticalAlignment = -4160
$q4.Range('A9').Borders.LineStyle = 1
$q4.Range('B9').NumberFormat = '@'
$q4.Range('B9').Value = "501209"
$q4.Range('C9').NumberFormat = '@'
$q4.Range('C9').Value = "银华富久食品饮料精选混合（LOF）A"
$q4.Range('D9').NumberFormat = '@'
$q4.Range('D9').Value = "6.30"
$q4.Range('E9').NumberFormat = '@'
$q4.Range('E9').Value = "90.18"
$q4.Range('F9').NumberFormat = '@'
$q4.Range('F9').Value = "4.40"
$q4.Range('G9').NumberFormat = '@'
$q4.Range('G9').Value = "0.2772"
$q4.Range('H9').Value = 9

# Row 10: fund 004237
$q4.Range('A10').Value = 8
$q4.Range('A10').Font.Bold = $true
$q4.Range('A10').HorizontalAlignment = -4108
$q4.Range('A10').VerticalAlignment = -4160
$q4.Range('A10').Borders.LineStyle = 1
$q4.Range('B10').NumberFormat = '@'
$q4.Range('B10').Value = "004237"
$q4.Range('C10').NumberFormat = '@'
$q4.Range('C10').Value = "中欧新蓝筹混合 -C"
$q4.Range('D10').NumberFormat = '@'
$q4.Range('D10').Value = "10.01"
$q4.Range('E10').NumberFormat = '@'
$q4.Range('E10').Value = "78.23"
$q4.Range('F10').NumberFormat = '@'
$q4.Range('F10').Value = "2.02"
$q4.Range('G10').NumberFormat = '@'
$q4.Range('G10').Value = "0.2022"
$q4.Range('H10').Value = 5

# Row 11: fund 121008
$q4.Range('A11').Value = 9
$q4.Range('A11').Font.Bold = $true
$q4.Range('A11').HorizontalAlignment = -4108
$q4.Range('A11').VerticalAlignment = -4160
$q4.Range('A11').Borders.LineStyle = 1
$q4.Range('B11').NumberFormat = '@'
$q4.Range('B11').Value = "121008"
$q4.Range('C11').NumberFormat = '@'
$q4.Range('C11').Value = "国投瑞银成长优选混合"
$q4.Range('D11').NumberFormat = '@'
$q4.Range('D11').Value = "5.75"
$q4.Range('E11').NumberFormat = '@'
$q4.Range('E11').Value = "88.41"
$q4.Range('F11').NumberFormat = '@'
$q4.Range('F11').Value = "3.10"
$q4.Range('G11').NumberFormat = '@'
$q4.Range('G11').Value = "0.1782"
$q4.Range('H11').Value = 10

# Row 12: fund 001218
$q4.Range('A12').Value = 10
$q4.Range('A12').Font.Bold = $true
$q4.Range('A12').HorizontalAlignment = -4108
$q4.Range('A12').VerticalAlignment = -4160
$q4.Range('A12').Borders.LineStyle = 1
$q4.Range('B12').NumberFormat = '@'
$q4.Range('B12').Value = "001218"
$q4.Range('C12').NumberFormat = '@'
$q4.Range('C12').Value = "国投瑞银精选收益灵活配置混合A"
$q4.Range('D12').NumberFormat = '@'
$q4.Range('D12').Value = "4.17"
$q4.Range('E12').NumberFormat = '@'
$q4.Range('E12').Value = "89.59"
$q4.Range('F12').NumberFormat = '@'
$q4.Range('F12').Value = "3.79"
$q4.Range('G12').NumberFormat = '@'
$q4.Range('G12').Value = "0.1580"
$q4.Range('H12').Value = 9

# Row 13: fund 015210
$q4.Range('A13').Value = 11
$q4.Range('A13').Font.Bold = $true
$q4.Range('A13').HorizontalAlignment = -4108
$q4.Range('A13').VerticalAlignment = -4160
$q4.Range('A13').Borders.LineStyle = 1
$q4.Range('B13').NumberFormat = '@'
$q4.Range('B13').Value = "015210"
$q4.Range('C13').NumberFormat = '@'
$q4.Range('C13').Value = "前海开源沪港深农业混合（LOF）C"
$q4.Range('D13').NumberFormat = '@'
$q4.Range('D13').Value = "1.94"
$q4.Range('E13').NumberFormat = '@'
$q4.Range('E13').Value = "88.37"
$q4.Range('F13').NumberFormat = '@'
$q4.Range('F13').Value = "7.77"
$q4.Range('G13').NumberFormat = '@'
$q4.Range('G13').Value = "0.1507"
$q4.Range('H13').Value = 1

# Row 14: fund 010418
$q4.Range('A14').Value = 12
$q4.Range('A14').Font.Bold = $true
$q4.Range('A14').HorizontalAlignment = -4108
$q4.Range('A14').VerticalAlignment = -4160
$q4.Range('A14').Borders.LineStyle = 1
$q4.Range('B14').NumberFormat = '@'
$q4.Range('B14').Value = "010418"
$q4.Range('C14').NumberFormat = '@'
$q4.Range('C14').Value = "财通景气行业混合A"
$q4.Range('D14').NumberFormat = '@'
$q4.Range('D14').Value = "2.60"
$q4.Range('E14').NumberFormat = '@'
$q4.Range('E14').Value = "86.37"
$q4.Range('F14').NumberFormat = '@'
$q4.Range('F14').Value = "5.45"
$q4.Range('G14').NumberFormat = '@'
$q4.Range('G14').Value = "0.1417"
$q4.Range('H14').Value = 8

# Row 15: fund 210004
$q4.Range('A15').Value = 13
$q4.Range('A15').Font.Bold = $true
$q4.Range('A15').HorizontalAlignment = -4108
$q4.Range('A15').VerticalAlignment = -4160
$q4.Range('A15').Borders.LineStyle = 1
$q4.Range('B15').NumberFormat = '@'
$q4.Range('B15').Value = "210004"
$q4.Range('C15').NumberFormat = '@'
$q4.Range('C15').Value = "金鹰稳健成长混合"
$q4.Range('D15').NumberFormat = '@'
$q4.Range('D15').Value = "5.14"
$q4.Range('E15').NumberFormat = '@'
$q4.Range('E15').Value = "94.39"
$q4.Range('F15').NumberFormat = '@'
$q4.Range('F15').Value = "2.64"
$q4.Range('G15').NumberFormat = '@'
$q4.Range('G15').Value = "0.1357"
$q4.Range('H15').Value = 7

# Row 16: fund 010423
$q4.Range('A16').Value = 14
$q4.Range('A16').Font.Bold = $true
$q4.Range('A16').HorizontalAlignment = -4108
$q4.Range('A16').VerticalAlignment = -4160
$q4.Range('A16').Borders.LineStyle = 1
$q4.Range('B16').NumberFormat = '@'
$q4.Range('B16').Value = "010423"
$q4.Range('C16').NumberFormat = '@'
$q4.Range('C16').Value = "国投瑞银价值成长一年持有期混合A"
$q4.Range('D16').NumberFormat = '@'
$q4.Range('D16').Value = "3.85"
$q4.Range('E16').NumberFormat = '@'
$q4.Range('E16').Value = "91.00"
$q4.Range('F16').NumberFormat = '@'
$q4.Range('F16').Value = "3.12"
$q4.Range('G16').NumberFormat = '@'
$q4.Range('G16').Value = "0.1201"
$q4.Range('H16').Value = 10

# Row 17: fund 501015
$q4.Range('A17').Value = 15
$q4.Range('A17').Font.Bold = $true
$q4.Range('A17').HorizontalAlignment = -4108
$q4.Range('A17').VerticalAlignment = -4160
$q4.Range('A17').Borders.LineStyle = 1
$q4.Range('B17').NumberFormat = '@'
$q4.Range('B17').Value = "501015"
$q4.Range('C17').NumberFormat = '@'
$q4.Range('C17').Value = "财通多策略升级混合（LOF）A"
$q4.Range('D17').NumberFormat = '@'
$q4.Range('D17').Value = "1.97"
$q4.Range('E17').NumberFormat = '@'
$q4.Range('E17').Value = "86.66"
$q4.Range('F17').NumberFormat = '@'
$q4.Range('F17').Value = "6.09"
$q4.Range('G17').NumberFormat = '@'
$q4.Range('G17').Value = "0.1200"
$q4.Range('H17').Value = 6

# Row 18: fund 501032
$q4.Range('A18').Value = 16
$q4.Range('A18').Font.Bold = $true
$q4.Range('A18').HorizontalAlignment = -4108
$q4.Range('A18').VerticalAlignment = -4160
$q4.Range('A18').Borders.LineStyle = 1
$q4.Range('B18').NumberFormat = '@'
$q4.Range('B18').Value = "501032"
$q4.Range('C18').NumberFormat = '@'
$q4.Range('C18').Value = "财通福盛多策略混合（LOF）A"
$q4.Range('D18').NumberFormat = '@'
$q4.Range('D18').Value = "2.16"
$q4.Range('E18').NumberFormat = '@'
$q4.Range('E18').Value = "87.95"
$q4.Range('F18').NumberFormat = '@'
$q4.Range('F18').Value = "5.41"
$q4.Range('G18').NumberFormat = '@'
$q4.Range('G18').Value = "0.1169"
$q4.Range('H18').Value = 3

# Row 19: fund 163503
$q4.Range('A19').Value = 17
$q4.Range('A19').Font.Bold = $true
$q4.Range('A19').HorizontalAlignment = -4108
$q4.Range('A19').VerticalAlignment = -4160
$q4.Range('A19').Borders.LineStyle = 1
$q4.Range('B19').NumberFormat = '@'
$q4.Range('B19').Value = "163503"
$q4.Range('C19').NumberFormat = '@'
$q4.Range('C19').Value = "天治核心成长混合（LOF）"
$q4.Range('D19').NumberFormat = '@'
$q4.Range('D19').Value = "3.14"
$q4.Range('E19').NumberFormat = '@'
$q4.Range('E19').Value = "93.97"
$q4.Range('F19').NumberFormat = '@'
$q4.Range('F19').Value = "3.42"
$q4.Range('G19').NumberFormat = '@'
$q4.Range('G19').Value = "0.1074"
$q4.Range('H19').Value = 2

# Row 20: fund 005959
$q4.Range('A20').Value = 18
$q4.Range('A20').Font.Bold = $true
$q4.Range('A20').HorizontalAlignment = -4108
$q4.Range('A20').VerticalAlignment = -4160
$q4.Range('A20').Borders.LineStyle = 1
$q4.Range('B20').NumberFormat = '@'
$q4.Range('B20').Value = "005959"
$q4.Range('C20').NumberFormat = '@'
$q4.Range('C20').Value = "财通新视野灵活配置混合C"
$q4.Range('D20').NumberFormat = '@'
$q4.Range('D20').Value = "1.61"
$q4.Range('E20').NumberFormat = '@'
$q4.Range('E20').Value = "86.51"
$q4.Range('F20').NumberFormat = '@'
$q4.Range('F20').Value = "6.50"
$q4.Range('G20').NumberFormat = '@'
$q4.Range('G20').Value = "0.1046"
$q4.Range('H20').Value = 6

# Row 21: fund 014064
$q4.Range('A21').Value = 19
$q4.Range('A21').Font.Bold = $true
$q4.Range('A21').HorizontalAlignment = -4108
$q4.Range('A21').VerticalAlignment = -4160
$q4.Range('A21').Borders.LineStyle = 1
$q4.Range('B21').NumberFormat = '@'
$q4.Range('B21').Value = "014064"
$q4.Range('C21').NumberFormat = '@'
$q4.Range('C21').Value = "银华农业产业股票C"
$q4.Range('D21').NumberFormat = '@'
$q4.Range('D21').Value = "1.12"
$q4.Range('E21').NumberFormat = '@'
$q4.Range('E21').Value = "93.03"
$q4.Range('F21').NumberFormat = '@'
$q4.Range('F21').Value = "5.33"
$q4.Range('G21').NumberFormat = '@'
$q4.Range('G21').Value = "0.0597"
$q4.Range('H21').Value = 7

# Row 22: fund 011181
$q4.Range('A22').Value = 20
$q4.Range('A22').Font.Bold = $true
$q4.Range('A22').HorizontalAlignment = -4108
$q4.Range('A22').VerticalAlignment = -4160
$q4.Range('A22').Borders.LineStyle = 1
$q4.Range('B22').NumberFormat = '@'
$q4.Range('B22').Value = "011181"
$q4.Range('C22').NumberFormat = '@'
$q4.Range('C22').Value = "长盛成长龙头混合A"
$q4.Range('D22').NumberFormat = '@'
$q4.Range('D22').Value = "1.03"
$q4.Range('E22').NumberFormat = '@'
$q4.Range('E22').Value = "87.87"
$q4.Range('F22').NumberFormat = '@'
$q4.Range('F22').Value = "4.26"
$q4.Range('G22').NumberFormat = '@'
$q4.Range('G22').Value = "0.0439"
$q4.Range('H22').Value = 8

# Row 23: fund 014210
$q4.Range('A23').Value = 21
$q4.Range('A23').Font.Bold = $true
$q4.Range('A23').HorizontalAlignment = -4108
$q4.Range('A23').VerticalAlignment = -4160
$q4.Range('A23').Borders.LineStyle = 1
$q4.Range('B23').NumberFormat = '@'
$q4.Range('B23').Value = "014210"
$q4.Range('C23').NumberFormat = '@'
$q4.Range('C23').Value = "国投瑞银竞争优势混合A"
$q4.Range('D23').NumberFormat = '@'
$q4.Range('D23').Value = "1.35"
$q4.Range('E23').NumberFormat = '@'
$q4.Range('E23').Value = "88.62"
$q4.Range('F23').NumberFormat = '@'
$q4.Range('F23').Value = "3.12"
$q4.Range('G23').NumberFormat = '@'
$q4.Range('G23').Value = "0.0421"
$q4.Range('H23').Value = 10

# Row 24: fund 005851
$q4.Range('A24').Value = 22
$q4.Range('A24').Font.Bold = $true
$q4.Range('A24').HorizontalAlignment = -4108
$q4.Range('A24').VerticalAlignment = -4160
$q4.Range('A24').Borders.LineStyle = 1
$q4.Range('B24').NumberFormat = '@'
$q4.Range('B24').Value = "005851"
$q4.Range('C24').NumberFormat = '@'
$q4.Range('C24').Value = "财通新视野灵活配置混合A"
$q4.Range('D24').NumberFormat = '@'
$q4.Range('D24').Value = "0.63"
$q4.Range('E24').NumberFormat = '@'
$q4.Range('E24').Value = "86.51"
$q4.Range('F24').NumberFormat = '@'
$q4.Range('F24').Value = "6.50"
$q4.Range('G24').NumberFormat = '@'
$q4.Range('G24').Value = "0.0410"
$q4.Range('H24').Value = 6

# Row 25: fund 013027
$q4.Range('A25').Value = 23
$q4.Range('A25').Font.Bold = $true
$q4.Range('A25').HorizontalAlignment = -4108
$q4.Range('A25').VerticalAlignment = -4160
$q4.Range('A25').Borders.LineStyle = 1
$q4.Range('B25').NumberFormat = '@'
$q4.Range('B25').Value = "013027"
$q4.Range('C25').NumberFormat = '@'
$q4.Range('C25').Value = "银华富久食品饮料精选混合（LOF）C"
$q4.Range('D25').NumberFormat = '@'
$q4.Range('D25').Value = "0.92"
$q4.Range('E25').NumberFormat = '@'
$q4.Range('E25').Value = "90.18"
$q4.Range('F25').NumberFormat = '@'
$q4.Range('F25').Value = "4.40"
$q4.Range('G25').NumberFormat = '@'
$q4.Range('G25').Value = "0.0405"
$q4.Range('H25').Value = 9

# Row 26: fund 015271
$q4.Range('A26').Value = 24
$q4.Range('A26').Font.Bold = $true
$q4.Range('A26').HorizontalAlignment = -4108
$q4.Range('A26').VerticalAlignment = -4160
$q4.Range('A26').Borders.LineStyle = 1
$q4.Range('B26').NumberFormat = '@'
$q4.Range('B26').Value = "015271"
$q4.Range('C26').NumberFormat = '@'
$q4.Range('C26').Value = "财通多策略升级混合（LOF）C"
$q4.Range('D26').NumberFormat = '@'
$q4.Range('D26').Value = "0.57"
$q4.Range('E26').NumberFormat = '@'
$q4.Range('E26').Value = "86.66"
$q4.Range('F26').NumberFormat = '@'
$q4.Range('F26').Value = "6.09"
$q4.Range('G26').NumberFormat = '@'
$q4.Range('G26').Value = "0.0347"
$q4.Range('H26').Value = 6

# Row 27: fund 001940
$q4.Range('A27').Value = 25
$q4.Range('A27').Font.Bold = $true
$q4.Range('A27').HorizontalAlignment = -4108
$q4.Range('A27').VerticalAlignment = -4160
$q4.Range('A27').Borders.LineStyle = 1
$q4.Range('B27').NumberFormat = '@'
$q4.Range('B27').Value = "001940"
$q4.Range('C27').NumberFormat = '@'
$q4.Range('C27').Value = "农银汇理现代农业加灵活配置混合"
$q4.Range('D27').NumberFormat = '@'
$q4.Range('D27').Value = "1.05"
$q4.Range('E27').NumberFormat = '@'
$q4.Range('E27').Value = "73.78"
$q4.Range('F27').NumberFormat = '@'
$q4.Range('F27').Value = "3.02"
$q4.Range('G27').NumberFormat = '@'
$q4.Range('G27').Value = "0.0317"
$q4.Range('H27').Value = 10

# Row 28: fund 010637
$q4.Range('A28').Value = 26
$q4.Range('A28').Font.Bold = $true
$q4.Range('A28').HorizontalAlignment = -4108
$q4.Range('A28').VerticalAlignment = -4160
$q4.Range('A28').Borders.LineStyle = 1
$q4.Range('B28').NumberFormat = '@'
$q4.Range('B28').Value = "010637"
$q4.Range('C28').NumberFormat = '@'
$q4.Range('C28').Value = "财通安盈混合C"
$q4.Range('D28').NumberFormat = '@'
$q4.Range('D28').Value = "1.58"
$q4.Range('E28').NumberFormat = '@'
$q4.Range('E28').Value = "36.44"
$q4.Range('F28').NumberFormat = '@'
$q4.Range('F28').Value = "1.91"
$q4.Range('G28').NumberFormat = '@'
$q4.Range('G28').Value = "0.0302"
$q4.Range('H28').Value = 5

# Row 29: fund 001885
$q4.Range('A29').Value = 27
$q4.Range('A29').Font.Bold = $true
$q4.Range('A29').HorizontalAlignment = -4108
$q4.Range('A29').VerticalAlignment = -4160
$q4.Range('A29').Borders.LineStyle = 1
$q4.Range('B29').NumberFormat = '@'
$q4.Range('B29').Value = "001885"
$q4.Range('C29').NumberFormat = '@'
$q4.Range('C29').Value = "中欧新蓝筹混合 -E"
$q4.Range('D29').NumberFormat = '@'
$q4.Range('D29').Value = "1.46"
$q4.Range('E29').NumberFormat = '@'
$q4.Range('E29').Value = "78.23"
$q4.Range('F29').NumberFormat = '@'
$q4.Range('F29').Value = "2.02"
$q4.Range('G29').NumberFormat = '@'
$q4.Range('G29').Value = "0.0295"
$q4.Range('H29').Value = 5

# Row 30: fund 000845
$q4.Range('A30').Value = 28
$q4.Range('A30').Font.Bold = $true
$q4.Range('A30').HorizontalAlignment = -4108
$q4.Range('A30').VerticalAlignment = -4160
$q4.Range('A30').Borders.LineStyle = 1
$q4.Range('B30').NumberFormat = '@'
$q4.Range('B30').Value = "000845"
$q4.Range('C30').NumberFormat = '@'
$q4.Range('C30').Value = "国投瑞银信息消费灵活配置混合"
$q4.Range('D30').NumberFormat = '@'
$q4.Range('D30').Value = "0.49"
$q4.Range('E30').NumberFormat = '@'
$q4.Range('E30').Value = "90.27"
$q4.Range('F30').NumberFormat = '@'
$q4.Range('F30').Value = "3.57"
$q4.Range('G30').NumberFormat = '@'
$q4.Range('G30').Value = "0.0175"
$q4.Range('H30').Value = 10

# Row 31: fund 010636
$q4.Range('A31').Value = 29
$q4.Range('A31').Font.Bold = $true
$q4.Range('A31').HorizontalAlignment = -4108
$q4.Range('A31').VerticalAlignment = -4160
$q4.Range('A31').Borders.LineStyle = 1
$q4.Range('B31').NumberFormat = '@'
$q4.Range('B31').Value = "010636"
$q4.Range('C31').NumberFormat = '@'
$q4.Range('C31').Value = "财通安盈混合A"
$q4.Range('D31').NumberFormat = '@'
$q4.Range('D31').Value = "0.88"
$q4.Range('E31').NumberFormat = '@'
$q4.Range('E31').Value = "36.44"
$q4.Range('F31').NumberFormat = '@'
$q4.Range('F31').Value = "1.91"
$q4.Range('G31').NumberFormat = '@'
$q4.Range('G31').Value = "0.0168"
$q4.Range('H31').Value = 5

# Row 32: fund 010188
$q4.Range('A32').Value = 30
$q4.Range('A32').Font.Bold = $true
$q4.Range('A32').HorizontalAlignment = -4108
$q4.Range('A32').VerticalAlignment = -4160
$q4.Range('A32').Borders.LineStyle = 1
$q4.Range('B32').NumberFormat = '@'
$q4.Range('B32').Value = "010188"
$q4.Range('C32').NumberFormat = '@'
$q4.Range('C32').Value = "中欧添益一年持有期混合A"
$q4.Range('D32').NumberFormat = '@'
$q4.Range('D32').Value = "1.47"
$q4.Range('E32').NumberFormat = '@'
$q4.Range('E32').Value = "27.97"
$q4.Range('F32').NumberFormat = '@'
$q4.Range('F32').Value = "1.06"
$q4.Range('G32').NumberFormat = '@'
$q4.Range('G32').Value = "0.0156"
$q4.Range('H32').Value = 10

# Row 33: fund 001520
$q4.Range('A33').Value = 31
$q4.Range('A33').Font.Bold = $true
$q4.Range('A33').HorizontalAlignment = -4108
$q4.Range('A33').VerticalAlignment = -4160
$q4.Range('A33').Borders.LineStyle = 1
$q4.Range('B33').NumberFormat = '@'
$q4.Range('B33').Value = "001520"
$q4.Range('C33').NumberFormat = '@'
$q4.Range('C33').Value = "国投瑞银研究精选股票"
$q4.Range('D33').NumberFormat = '@'
$q4.Range('D33').Value = "0.44"
$q4.Range('E33').NumberFormat = '@'
$q4.Range('E33').Value = "88.92"
$q4.Range('F33').NumberFormat = '@'
$q4.Range('F33').Value = "3.15"
$q4.Range('G33').NumberFormat = '@'
$q4.Range('G33').Value = "0.0139"
$q4.Range('H33').Value = 10

# Row 34: fund 010424
$q4.Range('A34').Value = 32
$q4.Range('A34').Font.Bold = $true
$q4.Range('A34').HorizontalAlignment = -4108
$q4.Range('A34').VerticalAlignment = -4160
$q4.Range('A34').Borders.LineStyle = 1
$q4.Range('B34').NumberFormat = '@'
$q4.Range('B34').Value = "010424"
$q4.Range('C34').NumberFormat = '@'
$q4.Range('C34').Value = "国投瑞银价值成长一年持有期混合C"
$q4.Range('D34').NumberFormat = '@'
$q4.Range('D34').Value = "0.29"
$q4.Range('E34').NumberFormat = '@'
$q4.Range('E34').Value = "91.00"
$q4.Range('F34').NumberFormat = '@'
$q4.Range('F34').Value = "3.12"
$q4.Range('G34').NumberFormat = '@'
$q4.Range('G34').Value = "0.0090"
$q4.Range('H34').Value = 10

# Row 35: fund 010189
$q4.Range('A35').Value = 33
$q4.Range('A35').Font.Bold = $true
$q4.Range('A35').HorizontalAlignment = -4108
$q4.Range('A35').VerticalAlignment = -4160
$q4.Range('A35').Borders.LineStyle = 1
$q4.Range('B35').NumberFormat = '@'
$q4.Range('B35').Value = "010189"
$q4.Range('C35').NumberFormat = '@'
$q4.Range('C35').Value = "中欧添益一年持有期混合C"
$q4.Range('D35').NumberFormat = '@'
$q4.Range('D35').Value = "0.83"
$q4.Range('E35').NumberFormat = '@'
$q4.Range('E35').Value = "27.97"
$q4.Range('F35').NumberFormat = '@'
$q4.Range('F35').Value = "1.06"
$q4.Range('G35').NumberFormat = '@'
$q4.Range('G35').Value = "0.0088"
$q4.Range('H35').Value = 10

# Row 36: fund 002005
$q4.Range('A36').Value = 34
$q4.Range('A36').Font.Bold = $true
$q4.Range('A36').HorizontalAlignment = -4108
$q4.Range('A36').VerticalAlignment = -4160
$q4.Range('A36').Borders.LineStyle = 1
$q4.Range('B36').NumberFormat = '@'
$q4.Range('B36').Value = "002005"
$q4.Range('C36').NumberFormat = '@'
$q4.Range('C36').Value = "工银新得利混合"
$q4.Range('D36').NumberFormat = '@'
$q4.Range('D36').Value = "0.49"
$q4.Range('E36').NumberFormat = '@'
$q4.Range('E36').Value = "26.23"
$q4.Range('F36').NumberFormat = '@'
$q4.Range('F36').Value = "1.15"
$q4.Range('G36').NumberFormat = '@'
$q4.Range('G36').Value = "0.0056"
$q4.Range('H36').Value = 8

# Row 37: fund 006458
$q4.Range('A37').Value = 35
$q4.Range('A37').Font.Bold = $true
$q4.Range('A37').HorizontalAlignment = -4108
$q4.Range('A37').VerticalAlignment = -4160
$q4.Range('A37').Borders.LineStyle = 1
$q4.Range('B37').NumberFormat = '@'
$q4.Range('B37').Value = "006458"
$q4.Range('C37').NumberFormat = '@'
$q4.Range('C37').Value = "平安估值优势灵活配置混合C"
$q4.Range('D37').NumberFormat = '@'
$q4.Range('D37').Value = "0.10"
$q4.Range('E37').NumberFormat = '@'
$q4.Range('E37').Value = "50.77"
$q4.Range('F37').NumberFormat = '@'
$q4.Range('F37').Value = "3.54"
$q4.Range('G37').NumberFormat = '@'
$q4.Range('G37').Value = "0.0035"
$q4.Range('H37').Value = 8

# Row 38: fund 014211
$q4.Range('A38').Value = 36
$q4.Range('A38').Font.Bold = $true
$q4.Range('A38').HorizontalAlignment = -4108
$q4.Range('A38').VerticalAlignment = -4160
$q4.Range('A38').Borders.LineStyle = 1
$q4.Range('B38').NumberFormat = '@'
$q4.Range('B38').Value = "014211"
$q4.Range('C38').NumberFormat = '@'
$q4.Range('C38').Value = "国投瑞银竞争优势混合C"
$q4.Range('D38').NumberFormat = '@'
$q4.Range('D38').Value = "0.04"
$q4.Range('E38').NumberFormat = '@'
$q4.Range('E38').Value = "88.62"
$q4.Range('F38').NumberFormat = '@'
$q4.Range('F38').Value = "3.12"
$q4.Range('G38').NumberFormat = '@'
$q4.Range('G38').Value = "0.0012"
$q4.Range('H38').Value = 10

# Row 39: fund 011182
$q4.Range('A39').Value = 37
$q4.Range('A39').Font.Bold = $true
$q4.Range('A39').HorizontalAlignment = -4108
$q4.Range('A39').VerticalAlignment = -4160
$q4.Range('A39').Borders.LineStyle = 1
$q4.Range('B39').NumberFormat = '@'
$q4.Range('B39').Value = "011182"
$q4.Range('C39').NumberFormat = '@'
$q4.Range('C39').Value = "长盛成长龙头混合C"
$q4.Range('D39').NumberFormat = '@'
$q4.Range('D39').Value = "0.02"
$q4.Range('E39').NumberFormat = '@'
$q4.Range('E39').Value = "87.87"
$q4.Range('F39').NumberFormat = '@'
$q4.Range('F39').Value = "4.26"
$q4.Range('G39').NumberFormat = '@'
$q4.Range('G39').Value = "0.0009"
$q4.Range('H39').Value = 8

# Row 40: fund 017679
$q4.Range('A40').Value = 38
$q4.Range('A40').Font.Bold = $true
$q4.Range('A40').HorizontalAlignment = -4108
$q4.Range('A40').VerticalAlignment = -4160
$q4.Range('A40').Borders.LineStyle = 1
$q4.Range('B40').NumberFormat = '@'
$q4.Range('B40').Value = "017679"
$q4.Range('C40').NumberFormat = '@'
$q4.Range('C40').Value = "国投瑞银精选收益灵活配置混合C"
$q4.Range('D40').NumberFormat = '@'
$q4.Range('D40').Value = "0.00"
$q4.Range('E40').NumberFormat = '@'
$q4.Range('E40').Value = "89.59"
$q4.Range('F40').NumberFormat = '@'
$q4.Range('F40').Value = "3.79"
$q4.Range('G40').Value = 0
$q4.Range('H40').Value = 9

# Row 41: fund 006457
$q4.Range('A41').Value = 39
$q4.Range('A41').Font.Bold = $true
$q4.Range('A41').HorizontalAlignment = -4108
$q4.Range('A41').VerticalAlignment = -4160
$q4.Range('A41').Borders.LineStyle = 1
$q4.Range('B41').NumberFormat = '@'
$q4.Range('B41').Value = "006457"
$q4.Range('C41').NumberFormat = '@'
$q4.Range('C41').Value = "平安估值优势灵活配置混合A"
$q4.Range('D41').NumberFormat = '@'
$q4.Range('D41').Value = "0.00"
$q4.Range('E41').NumberFormat = '@'
$q4.Range('E41').Value = "50.77"
$q4.Range('F41').NumberFormat = '@'
$q4.Range('F41').Value = "3.54"
$q4.Range('G41').Value = 0
$q4.Range('H41').Value = 8

# Row 42: fund 014628
$q4.Range('A42').Value = 40
$q4.Range('A42').Font.Bold = $true
$q4.Range('A42').HorizontalAlignment = -4108
$q4.Range('A42').VerticalAlignment = -4160
$q4.Range('A42').Borders.LineStyle = 1
$q4.Range('B42').NumberFormat = '@'
$q4.Range('B42').Value = "014628"
$q4.Range('C42').NumberFormat = '@'
$q4.Range('C42').Value = "财通福盛多策略混合（LOF）C"
$q4.Range('D42').NumberFormat = '@'
$q4.Range('D42').Value = "0.00"
$q4.Range('E42').NumberFormat = '@'
$q4.Range('E42').Value = "87.95"
$q4.Range('F42').NumberFormat = '@'
$q4.Range('F42').Value = "5.41"
$q4.Range('G42').Value = 0
$q4.Range('H42').Value = 3

# Row 43: fund 016234
$q4.Range('A43').Value = 41
$q4.Range('A43').Font.Bold = $true
$q4.Range('A43').HorizontalAlignment = -4108
$q4.Range('A43').VerticalAlignment = -4160
$q4.Range('A43').Borders.LineStyle = 1
$q4.Range('B43').NumberFormat = '@'
$q4.Range('B43').Value = "016234"
$q4.Range('C43').NumberFormat = '@'
$q4.Range('C43').Value = "财通景气行业混合C"
$q4.Range('D43').NumberFormat = '@'
$q4.Range('D43').Value = "0.00"
$q4.Range('E43').NumberFormat = '@'
$q4.Range('E43').Value = "86.37"
$q4.Range('F43').NumberFormat = '@'
$q4.Range('F43').Value = "5.45"
$q4.Range('G43').Value = 0
$q4.Range('H43').Value = 8

